# Apply September Deskcount updates:
# - Greenwood Village (row 16), Tampa (row 38), Santiago (row 47), Sao Paulo (row 48)
#   "Include in Occupancy Calculation" flipped from Yes -> No
# - Melbourne (row 44) Deskcount updated 30 -> 32
# - Selection/view scrolled down to row 31, active cell C45

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F16").Value = "No"
$ws.Range("F38").Value = "No"
$ws.Range("C44").Value = 32
$ws.Range("F47").Value = "No"
$ws.Range("F48").Value = "No"

$ws.Range("C45").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 2
